$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Nacionalidad -> "Extranjera", Descuento -> "Si"
$ws.Range("I2").Value = "Extranjera"
$ws.Range("J2").Value = "Si"

# Row 3: Tipo -> "Otros", Descuento -> "No"
$ws.Range("H3").Value = "Otros"
$ws.Range("J3").Value = "No"

# New row 4: single-space value in Descuento column
$ws.Range("J4").Value = " "

# Update the active selection to match the saved state
$ws.Range("E9").Select()
